# Update cryptocurrency price/volume table (columns D and E) per the latest
# scrape of coinranking.com data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.575.43"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "3.364.26"
$ws.Range("E3").Value = "  -4.36%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -2.73%  "
$ws.Range("D8").Value = "3.356.43"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.627"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "3.906.57"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "3.378.64"
$ws.Range("E17").Value = "  -4.19%  "
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "64.566.21"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.980"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "576.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.49%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -7.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -5.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.368"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").Value = "3.114.42"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.11%  "
$ws.Range("E45").Value = "  -3.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0410"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.46%  "
